$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Transacciones": append two new transaction rows (176 and 177)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Transacciones")

# Row 176 - Salario Quincena (income)
$ws1.Range("A176").Value = 43600
$ws1.Range("B176").Value = 5827
$ws1.Range("C176").Value = "Salario Quincena"
$ws1.Range("D176").Value = "Salario"
$ws1.Range("E176").Value = "Ingreso"
$ws1.Range("F176").Value = "Tarjeta Santander"
$ws1.Range("G176").Value = "Optimen"
$ws1.Range("K176").Value = 4856.18
$ws1.Range("L176").Formula = "=L175+B176"
$ws1.Range("M176").Value = 5
$ws1.Range("N176").Formula = "=SUM(K176:M176)"
$ws1.Range("O176").Formula = "=N176-4000"
$ws1.Range("P176").Formula = "=O176-Ahorros!`$E`$4"

# Row 177 - Pago de Deudas (expense)
$ws1.Range("A177").Value = 43600
$ws1.Range("B177").Value = 2000
$ws1.Range("C177").Value = "Pago de Deudas"
$ws1.Range("D177").Value = "Pagos"
$ws1.Range("E177").Value = "Gasto"
$ws1.Range("F177").Value = "Tarjeta Santander"
$ws1.Range("G177").Value = "NA"
$ws1.Range("K177").Value = 4856.18
$ws1.Range("L177").Formula = "=L176-B177"
$ws1.Range("M177").Value = 5
$ws1.Range("N177").Formula = "=SUM(K177:M177)"
$ws1.Range("O177").Formula = "=N177-4000"
$ws1.Range("P177").Formula = "=O177-Ahorros!`$E`$4"

# Match the date-number formatting (s="1") and the "Bueno" highlight style
# (s="25") used by the rest of the table, by copying formats from the row
# directly above / the existing P column.
$ws1.Range("A175").Copy()
$ws1.Range("A176:A177").PasteSpecial(-4122)
$ws1.Range("P176:P177").Style = "Bueno"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "Deudas": append a new debt-payment row (15)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Deudas")

$ws2.Range("I15").Value = 43600
$ws2.Range("J15").Value = 2000
$ws2.Range("K15").Value = "Pago"

$ws2.Range("I14").Copy()
$ws2.Range("I15").PasteSpecial(-4122)
$ws2.Range("I15").Value = 43600

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "Ahorros": record the new 43600 savings deposit (row 15) and the
# related entry in H6/I6
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Ahorros")

$ws3.Range("H6").Value = 6000
$ws3.Range("I6").Value = 43600

$ws3.Range("A14").Copy()
$ws3.Range("A15").PasteSpecial(-4122)
$ws3.Range("A15").Value = 43600
$ws3.Range("B15").Value = 500

$ws3.Range("I5").Copy()
$ws3.Range("I6").PasteSpecial(-4122)
$ws3.Range("I6").Value = 43600

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Recalculate and restore the view/selection state shown by the diff:
#   Deudas   -> L15 selected
#   Ahorros  -> C15 selected
#   Transacciones -> Q177 selected, and left as the active sheet/tab
# ---------------------------------------------------------------------------
$excel.Calculate()

$ws2.Activate()
$ws2.Range("L15").Select()

$ws3.Activate()
$ws3.Range("C15").Select()

$ws1.Activate()
$ws1.Range("Q177").Select()
